# Update LR-pair data with newly computed TPM-based values.
# The sending/target cluster combinations that pair with "ECs" as the
# Target cluster are dropped entirely, and the remaining six combinations
# (rows 2-7) are rewritten with their freshly computed statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Icam4","Itgam","FAPs",2,0.6666666666666666,0.2564746666666666,0.7694239999999999,0.1818007399394835,0.1818007399394835,1,0.3333333333333333,0.005673666666666667,0.017021,0.1234828534325781,0.1234828534325781,0.001455151767111111,0.013096365904,0.02244927412388149,0.0224492741238815),
    @("ECs","Icam4","Itgam","MuSCs",2,0.6666666666666666,0.2564746666666666,0.7694239999999999,0.1818007399394835,0.1818007399394835,1,0.3333333333333333,0.04027333333333333,0.12082,0.8765171465674219,0.876517146567422,0.01032908974222222,0.09296180767999998,0.159351465815602,0.1593514658156021),
    @("FAPs","Icam4","Itgam","FAPs",3,1,1.085481666666667,3.256445,0.7694380609030022,0.7694380609030022,1,0.3333333333333333,0.005673666666666667,0.017021,0.1234828534325781,0.1234828534325781,0.006158661149444445,0.055427950345,0.09501240729993253,0.09501240729993254),
    @("FAPs","Icam4","Itgam","MuSCs",3,1,1.085481666666667,3.256445,0.7694380609030022,0.7694380609030022,1,0.3333333333333333,0.04027333333333333,0.12082,0.8765171465674219,0.876517146567422,0.04371596498888889,0.3934436849,0.6744256536030696,0.6744256536030697),
    @("MuSCs","Icam4","Itgam","FAPs",2,0.6666666666666666,0.06878966666666667,0.206369,0.0487611991575143,0.0487611991575143,1,0.3333333333333333,0.005673666666666667,0.017021,0.1234828534325781,0.1234828534325781,0.0003902896387777778,0.003512606749,0.006021172008764089,0.00602117200876409),
    @("MuSCs","Icam4","Itgam","MuSCs",2,0.6666666666666666,0.06878966666666667,0.206369,0.0487611991575143,0.0487611991575143,1,0.3333333333333333,0.04027333333333333,0.12082,0.8765171465674219,0.876517146567422,0.002770389175555556,0.02493350258,0.04274002714875021,0.04274002714875021)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($i + 2, $j + 1).Value = $row[$j]
    }
}

# Rows 8-10 (the old "Target cluster = ECs" combos for MuSCs-sending rows,
# plus the now-unused trailing rows) are removed so the sheet shrinks to
# A1:T7.
$ws.Range("A8:T10").Delete() | Out-Null
